# Fruta / hortaliza, semanal
# The sheet's 40 data records (rows 2-41) get re-shuffled: each record's
# Fecha (D), Calidad (L), Volumen (M), Precio minimo/maximo/promedio
# (N/O/P), Origen (R) and Precio $/Kg (S) move to a different row, while
# all other columns (A,B,C,E,F,G,H,I,J,K,Q,T) stay put on their row.
#
# destination row -> source row (i.e. row <key> ends up holding the data
# that used to live on row <value>)
$map = @{
    2  = 18
    3  = 19
    4  = 34
    5  = 35
    6  = 33
    7  = 4
    8  = 5
    9  = 23
    10 = 14
    11 = 15
    12 = 39
    13 = 40
    14 = 38
    15 = 36
    16 = 27
    17 = 24
    18 = 30
    19 = 31
    20 = 25
    21 = 9
    22 = 28
    23 = 29
    24 = 12
    25 = 13
    26 = 32
    27 = 20
    28 = 21
    29 = 11
    30 = 10
    31 = 22
    32 = 26
    33 = 37
    34 = 6
    35 = 7
    36 = 8
    37 = 16
    38 = 17
    39 = 41
    40 = 2
    41 = 3
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot every affected cell's current value BEFORE writing anything,
# since this is a full permutation (not a simple shift) and rows read
# from and written to overlap.
$snapshot = @{}
for ($row = 2; $row -le 41; $row++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Write each destination row using the snapshotted source row's values.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
